$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.316.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.871.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4709"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2885"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("E9").Value = "  +1.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08049"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.876.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.145"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6872"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "271.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.311.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("E18").Value = "  +5.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007762"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.26%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.117.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.323"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.03%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.219"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.360"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.958"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09935"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.366"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.467"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.082"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04708"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("E35").Value = "  +0.91%  "
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.707"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01885"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.655"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.302"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.961"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8433"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4171"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.0000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.307"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.091"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "933.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05683"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.65%  "